# 自动更新Excel文件 - 2025-11-26 23:12:23
# For every data row (row 2..last used row) on the active sheet:
#   - Column D = 总天 (total days for the cycle)
#   - Column E = 剩余 (days remaining)
#   - Column F = 开始时间 (cycle start date, stored as an 8-digit yyyyMMdd integer)
# Each day the remaining-days counter ticks down by 1. When it would reach 0
# (i.e. the stored value is 1), the cycle restarts: remaining resets to the
# full cycle length (D) and the start date is pushed forward by D days.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($dVal -eq $null -or $eVal -eq $null -or $fVal -eq $null) {
        continue
    }

    $fText = [string]([int64]$fVal)

    # Only rows whose start date is a well-formed 8-digit yyyyMMdd value are
    # advanced; malformed dates are left untouched (mirrors source data).
    if ($fText.Length -ne 8) {
        continue
    }

    $year = [int]$fText.Substring(0, 4)
    $month = [int]$fText.Substring(4, 2)
    $day = [int]$fText.Substring(6, 2)

    $totalDays = [int]$dVal
    $remaining = [int]$eVal

    if ($remaining -gt 1) {
        $eCell.Value2 = $remaining - 1
    }
    else {
        $startDate = Get-Date -Year $year -Month $month -Day $day
        $newStartDate = $startDate.AddDays($totalDays)

        $eCell.Value2 = $totalDays
        $fCell.Value2 = [int64]($newStartDate.ToString("yyyyMMdd"))
    }
}
